$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Efna3/Epha3 -> FAPs
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.046374
$ws.Range("H2").Value = 0.139122
$ws.Range("I2").Value = 0.6592866045237633
$ws.Range("J2").Value = 0.6592866045237632
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.097779
$ws.Range("N2").Value = 78.29333700000001
$ws.Range("O2").Value = 0.9922055808976035
$ws.Range("P2").Value = 0.9922055808976036
$ws.Range("Q2").Value = 1.210258403346
$ws.Range("R2").Value = 10.892325630114
$ws.Range("S2").Value = 0.6541478484195091
$ws.Range("T2").Value = 0.6541478484195091

# Row 3: ECs -> Efna3/Epha3 -> MuSCs
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.046374
$ws.Range("H3").Value = 0.139122
$ws.Range("I3").Value = 0.6592866045237633
$ws.Range("J3").Value = 0.6592866045237632
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.205015
$ws.Range("N3").Value = 0.6150450000000001
$ws.Range("O3").Value = 0.007794419102396499
$ws.Range("P3").Value = 0.007794419102396499
$ws.Range("Q3").Value = 0.00950736561
$ws.Range("R3").Value = 0.08556629049
$ws.Range("S3").Value = 0.005138756104254147
$ws.Range("T3").Value = 0.005138756104254146

# Row 4: MuSCs -> Efna3/Epha3 -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.02396566666666667
$ws.Range("H4").Value = 0.071897
$ws.Range("I4").Value = 0.3407133954762367
$ws.Range("J4").Value = 0.3407133954762367
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.097779
$ws.Range("N4").Value = 78.29333700000001
$ws.Range("O4").Value = 0.9922055808976035
$ws.Range("P4").Value = 0.9922055808976036
$ws.Range("Q4").Value = 0.6254506722543334
$ws.Range("R4").Value = 5.629056050289001
$ws.Range("S4").Value = 0.3380577324780944
$ws.Range("T4").Value = 0.3380577324780944

# Row 5 (new): MuSCs -> Efna3/Epha3 -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02396566666666667
$ws.Range("H5").Value = 0.071897
$ws.Range("I5").Value = 0.3407133954762367
$ws.Range("J5").Value = 0.3407133954762367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.205015
$ws.Range("N5").Value = 0.6150450000000001
$ws.Range("O5").Value = 0.007794419102396499
$ws.Range("P5").Value = 0.007794419102396499
$ws.Range("Q5").Value = 0.004913321151666667
$ws.Range("R5").Value = 0.04421989036500001
$ws.Range("S5").Value = 0.002655662998142353
$ws.Range("T5").Value = 0.002655662998142353
